{"js": "// Add the institutional affiliation line right under the \"Edison Achalma\"\n// byline (the paragraph styled \"Author\" that sits directly below the title),\n// as a new paragraph that reuses the same \"Author\" style.\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text,items/style\");\nawait context.sync();\n\n// Locate the byline paragraph by its exact text + style, instead of a\n// hard-coded index, since \"Edison Achalma\" also appears later in the\n// document (author note, CRediT statement, correspondence line).\nlet targetIndex = -1;\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  const para = paragraphs.items[i];\n  if (para.text === \"Edison Achalma\" && para.style === \"Author\") {\n    targetIndex = i;\n    break;\n  }\n}\nif (targetIndex === -1) {\n  throw new Error('Could not find the \"Edison Achalma\" byline paragraph.');\n}\n\n// Insert the new text + a paragraph break right before the paragraph that\n// follows the byline. Doing it this way (rather than calling\n// `insertParagraph` directly on the byline paragraph) leaves the existing\n// \"Edison Achalma\" paragraph completely untouched.\nconst followingPara = paragraphs.items[targetIndex + 1];\nconst insertionRange = followingPara.getRange(Word.RangeLocation.start);\ninsertionRange.insertText(\n  \"Escuela Profesional de Econom\u00eda, Universidad Nacional de San Crist\u00f3bal de Huamanga\\n\",\n  Word.InsertLocation.before\n);\nawait context.sync();\n\n// The newly split-off paragraph inherited the following paragraph's style;\n// give it the \"Author\" style to match the byline block.\nconst refreshedParagraphs = context.document.body.paragraphs;\nrefreshedParagraphs.load(\"items/text\");\nawait context.sync();\nconst newPara = refreshedParagraphs.items[targetIndex + 1];\nnewPara.style = \"Author\";\nawait context.sync();\n", "ps1": "# Add the institutional affiliation line right under the \"Edison Achalma\"\n# byline (the paragraph styled \"Author\" that sits directly below the title),\n# as a new paragraph that reuses the same \"Author\" style.\n\n$d = $word.ActiveDocument\n\n# Locate the byline paragraph by its exact text + style, instead of a\n# hard-coded index, since \"Edison Achalma\" also appears later in the\n# document (author note, CRediT statement, correspondence line).\n$targetPara = $null\nforeach ($p in $d.Paragraphs) {\n    $text = $p.Range.Text -replace \"[\\r\\a]+$\", \"\"\n    if ($text -eq \"Edison Achalma\" -and $p.Style.NameLocal -eq \"Author\") {\n        $targetPara = $p\n        break\n    }\n}\nif ($null -eq $targetPara) {\n    throw \"Could not find the 'Edison Achalma' byline paragraph\"\n}\n\n# Insert a new paragraph mark right after the byline, then fill it in. This\n# leaves the existing \"Edison Achalma\" paragraph completely untouched.\n$insertionPoint = $d.Range($targetPara.Range.End, $targetPara.Range.End)\n$insertionPoint.InsertParagraphAfter()\n\n$newPara = $targetPara.Next()\n$newPara.Style = \"Author\"\n$newPara.Range.Text = \"Escuela Profesional de Econom\u00eda, Universidad Nacional de San Crist\u00f3bal de Huamanga\"\n"}
